$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "Haggard'OlGoblin"
$ws.Range("C35").Value = "Haggard 'Ol Goblin"
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = "Labyrinth"
$ws.Range("I35").Value = 16
$ws.Range("J35").Value = 148
